$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44685
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 1500
$ws.Cells.Item(2, 12).Value = 2000
$ws.Cells.Item(2, 13).Value = 1750
$ws.Cells.Item(2, 16).Value = 583
$ws.Cells.Item(3, 4).Value = 45041
$ws.Cells.Item(3, 10).Value = 1160
$ws.Cells.Item(3, 11).Value = 2000
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = 2250
$ws.Cells.Item(3, 16).Value = 750
$ws.Cells.Item(4, 4).Value = 45077
$ws.Cells.Item(4, 10).Value = 760
$ws.Cells.Item(4, 11).Value = 2000
$ws.Cells.Item(4, 12).Value = 2500
$ws.Cells.Item(4, 13).Value = 2250
$ws.Cells.Item(4, 16).Value = 750
$ws.Cells.Item(5, 4).Value = 45034
$ws.Cells.Item(5, 10).Value = 1100
$ws.Cells.Item(5, 11).Value = 2000
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2250
$ws.Cells.Item(5, 16).Value = 750
$ws.Cells.Item(6, 4).Value = 45028
$ws.Cells.Item(6, 10).Value = 1000
$ws.Cells.Item(6, 11).Value = 2000
$ws.Cells.Item(6, 12).Value = 2500
$ws.Cells.Item(6, 13).Value = 2250
$ws.Cells.Item(6, 16).Value = 750
$ws.Cells.Item(7, 4).Value = 45006
$ws.Cells.Item(7, 10).Value = 1100
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 12).Value = 2500
$ws.Cells.Item(7, 13).Value = 2250
$ws.Cells.Item(7, 16).Value = 750
$ws.Cells.Item(8, 4).Value = 44985
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 2000
$ws.Cells.Item(8, 12).Value = 2500
$ws.Cells.Item(8, 13).Value = 2250
$ws.Cells.Item(8, 16).Value = 750
$ws.Cells.Item(9, 4).Value = 45013
$ws.Cells.Item(9, 10).Value = 1100
$ws.Cells.Item(9, 11).Value = 2000
$ws.Cells.Item(9, 12).Value = 2500
$ws.Cells.Item(9, 13).Value = 2250
$ws.Cells.Item(9, 16).Value = 750
$ws.Cells.Item(10, 4).Value = 45084
$ws.Cells.Item(10, 10).Value = 900
$ws.Cells.Item(10, 11).Value = 2000
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2250
$ws.Cells.Item(10, 16).Value = 750
$ws.Cells.Item(11, 4).Value = 44978
$ws.Cells.Item(11, 10).Value = 1000
$ws.Cells.Item(11, 11).Value = 1800
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 1900
$ws.Cells.Item(11, 16).Value = 633
$ws.Cells.Item(12, 4).Value = 44911
$ws.Cells.Item(12, 10).Value = 700
$ws.Cells.Item(12, 11).Value = 1800
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 1900
$ws.Cells.Item(12, 16).Value = 633
$ws.Cells.Item(13, 4).Value = 44999
$ws.Cells.Item(13, 10).Value = 1100
$ws.Cells.Item(13, 11).Value = 2000
$ws.Cells.Item(13, 12).Value = 2500
$ws.Cells.Item(13, 13).Value = 2250
$ws.Cells.Item(13, 16).Value = 750
$ws.Cells.Item(14, 4).Value = 45007
$ws.Cells.Item(14, 10).Value = 1160
$ws.Cells.Item(14, 11).Value = 2000
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 13).Value = 2250
$ws.Cells.Item(14, 16).Value = 750
$ws.Cells.Item(15, 4).Value = 44971
$ws.Cells.Item(15, 10).Value = 1000
$ws.Cells.Item(15, 11).Value = 2000
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2250
$ws.Cells.Item(15, 16).Value = 750
$ws.Cells.Item(16, 4).Value = 45020
$ws.Cells.Item(16, 10).Value = 1200
$ws.Cells.Item(16, 11).Value = 2000
$ws.Cells.Item(16, 12).Value = 2500
$ws.Cells.Item(16, 13).Value = 2250
$ws.Cells.Item(16, 16).Value = 750
$ws.Cells.Item(17, 4).Value = 45035
$ws.Cells.Item(17, 10).Value = 1100
$ws.Cells.Item(17, 11).Value = 2000
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = 2250
$ws.Cells.Item(17, 16).Value = 750
$ws.Cells.Item(18, 4).Value = 44953
$ws.Cells.Item(18, 10).Value = 1000
$ws.Cells.Item(18, 11).Value = 2000
$ws.Cells.Item(18, 12).Value = 2500
$ws.Cells.Item(18, 13).Value = 2250
$ws.Cells.Item(18, 16).Value = 750
$ws.Cells.Item(19, 4).Value = 44848
$ws.Cells.Item(19, 10).Value = 1000
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 1750
$ws.Cells.Item(19, 16).Value = 583
$ws.Cells.Item(20, 4).Value = 45070
$ws.Cells.Item(20, 10).Value = 800
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2250
$ws.Cells.Item(20, 16).Value = 750
$ws.Cells.Item(21, 4).Value = 45062
$ws.Cells.Item(21, 10).Value = 1100
$ws.Cells.Item(21, 11).Value = 2000
$ws.Cells.Item(21, 12).Value = 2500
$ws.Cells.Item(21, 13).Value = 2250
$ws.Cells.Item(21, 16).Value = 750
$ws.Cells.Item(22, 4).Value = 44951
$ws.Cells.Item(22, 10).Value = 800
$ws.Cells.Item(22, 11).Value = 2000
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2250
$ws.Cells.Item(22, 16).Value = 750
$ws.Cells.Item(23, 4).Value = 44883
$ws.Cells.Item(23, 10).Value = 500
$ws.Cells.Item(23, 11).Value = 1800
$ws.Cells.Item(23, 12).Value = 2000
$ws.Cells.Item(23, 13).Value = 1900
$ws.Cells.Item(23, 16).Value = 633
$ws.Cells.Item(24, 4).Value = 44964
$ws.Cells.Item(24, 10).Value = 1000
$ws.Cells.Item(24, 11).Value = 2000
$ws.Cells.Item(24, 12).Value = 2500
$ws.Cells.Item(24, 13).Value = 2250
$ws.Cells.Item(24, 16).Value = 750
$ws.Cells.Item(25, 4).Value = 44827
$ws.Cells.Item(25, 10).Value = 1200
$ws.Cells.Item(25, 11).Value = 2000
$ws.Cells.Item(25, 12).Value = 2500
$ws.Cells.Item(25, 13).Value = 2250
$ws.Cells.Item(25, 16).Value = 750
$ws.Cells.Item(26, 4).Value = 44910
$ws.Cells.Item(26, 10).Value = 1000
$ws.Cells.Item(26, 11).Value = 1800
$ws.Cells.Item(26, 12).Value = 2000
$ws.Cells.Item(26, 13).Value = 1900
$ws.Cells.Item(26, 16).Value = 633
$ws.Cells.Item(27, 4).Value = 44970
$ws.Cells.Item(27, 10).Value = 800
$ws.Cells.Item(27, 11).Value = 2000
$ws.Cells.Item(27, 12).Value = 2500
$ws.Cells.Item(27, 13).Value = 2250
$ws.Cells.Item(27, 16).Value = 750
$ws.Cells.Item(28, 4).Value = 44992
$ws.Cells.Item(28, 10).Value = 1040
$ws.Cells.Item(28, 11).Value = 2000
$ws.Cells.Item(28, 12).Value = 2500
$ws.Cells.Item(28, 13).Value = 2250
$ws.Cells.Item(28, 16).Value = 750
$ws.Cells.Item(29, 4).Value = 44965
$ws.Cells.Item(29, 10).Value = 1120
$ws.Cells.Item(29, 11).Value = 2000
$ws.Cells.Item(29, 12).Value = 2500
$ws.Cells.Item(29, 13).Value = 2250
$ws.Cells.Item(29, 16).Value = 750
$ws.Cells.Item(30, 4).Value = 44881
$ws.Cells.Item(30, 10).Value = 500
$ws.Cells.Item(30, 11).Value = 1900
$ws.Cells.Item(30, 12).Value = 2000
$ws.Cells.Item(30, 13).Value = 1950
$ws.Cells.Item(30, 16).Value = 650
